# Update simulation result values (rows 2-29) with the latest run output
# from multiple_trains_crane_2_hostler_3_simulation_results.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.33349672386051632
$ws.Range("E2").Value = 0.4334967238605163
$ws.Range("K2").Value = 0.53349672386051628
$ws.Range("L2").Value = 0.56748017527631467
$ws.Range("N2").Value = 0.13349672386051631
$ws.Range("O2").Value = 0.46748017527631469
$ws.Range("D3").Value = 3.6676152300564349
$ws.Range("E3").Value = 3.7676152300564349
$ws.Range("I3").Value = 3.434199106974714
$ws.Range("K3").Value = 3.867615230056435
$ws.Range("L3").Value = 3.904873635092819
$ws.Range("N3").Value = 3.4676152300564351
$ws.Range("O3").Value = 3.50487363509282
$ws.Range("D4").Value = 3.5003456208799029
$ws.Range("E4").Value = 3.600345620879903
$ws.Range("I4").Value = 3.2668728499564779
$ws.Range("K4").Value = 3.700345620879903
$ws.Range("L4").Value = 3.744449172113133
$ws.Range("M4").Value = 3.3003456208799031
$ws.Range("O4").Value = 3.4444491721131332
$ws.Range("D5").Value = 0.33340955013056989
$ws.Range("E5").Value = 0.43340955013056992
$ws.Range("J5").Value = 0.13340955013056999
$ws.Range("K5").Value = 0.5334095501305699
$ws.Range("L5").Value = 0.59661577771035224
$ws.Range("O5").Value = 0.49661577771035231
$ws.Range("D6").Value = 2.333514031909476
$ws.Range("E6").Value = 2.4335140319094761
$ws.Range("I6").Value = 2.0668957364016158
$ws.Range("K6").Value = 2.5335140319094762
$ws.Range("L6").Value = 2.5805384368056341
$ws.Range("N6").Value = 2.133514031909475
$ws.Range("O6").Value = 2.480538436805634
$ws.Range("D7").Value = 3.5338451397709232
$ws.Range("E7").Value = 3.6338451397709228
$ws.Range("I7").Value = 3.2668728499564779
$ws.Range("K7").Value = 3.7338451397709229
$ws.Range("L7").Value = 3.794802905929739
$ws.Range("M7").Value = 3.333845139770923
$ws.Range("O7").Value = 3.4948029059297392
$ws.Range("D8").Value = 0.33336779600325711
$ws.Range("E8").Value = 0.43336779600325698
$ws.Range("K8").Value = 0.53336779600325701
$ws.Range("L8").Value = 0.58599516895840331
$ws.Range("M8").Value = 0.13336779600325699
$ws.Range("O8").Value = 0.48599516895840328
$ws.Range("D9").Value = 2.36691466396173
$ws.Range("E9").Value = 2.466914663961731
$ws.Range("I9").Value = 2.0668957364016158
$ws.Range("K9").Value = 2.5669146639617311
$ws.Range("L9").Value = 2.6313967473001849
$ws.Range("N9").Value = 2.1669146639617298
$ws.Range("O9").Value = 2.5313967473001848
$ws.Range("D10").Value = 0.73336779600325697
$ws.Range("E10").Value = 0.83336779600325694
$ws.Range("K10").Value = 0.93336779600325692
$ws.Range("L10").Value = 0.993915627385619
$ws.Range("M10").Value = 0.16675615851641629
$ws.Range("O10").Value = 0.89391562738561903
$ws.Range("D11").Value = 2.433409550130571
$ws.Range("E11").Value = 2.5334095501305711
$ws.Range("I11").Value = 2.0668957364016158
$ws.Range("K11").Value = 2.6334095501305712
$ws.Range("L11").Value = 2.6698999822299578
$ws.Range("N11").Value = 2.2003855938919421
$ws.Range("O11").Value = 2.5698999822299582
$ws.Range("D12").Value = 2.7669146639617308
$ws.Range("E12").Value = 2.8669146639617309
$ws.Range("I12").Value = 2.0668957364016158
$ws.Range("K12").Value = 2.966914663961731
$ws.Range("L12").Value = 3.0129026715314371
$ws.Range("N12").Value = 2.2338183854240818
$ws.Range("O12").Value = 2.912902671531437
$ws.Range("D13").Value = 2.8334095501305709
$ws.Range("E13").Value = 2.933409550130571
$ws.Range("I13").Value = 2.0668957364016158
$ws.Range("K13").Value = 3.033409550130572
$ws.Range("L13").Value = 3.0905323040723589
$ws.Range("N13").Value = 2.2672843189960399
$ws.Range("O13").Value = 2.9905323040723588
$ws.Range("B14").Value = 0.049144239375912897
$ws.Range("C14").Value = 0.049144239375912897
$ws.Range("D14").Value = 1.73340955013057
$ws.Range("E14").Value = 1.933409550130571
$ws.Range("F14").Value = 2.1001718484605831
$ws.Range("G14").Value = 2.0668957364016158
$ws.Range("O14").Value = 2.0510276090846702
$ws.Range("B15").Value = 0.043553692430931087
$ws.Range("C15").Value = 0.043553692430931087
$ws.Range("D15").Value = 1.133496723860516
$ws.Range("E15").Value = 1.333496723860516
$ws.Range("F15").Value = 4.0338124751583351
$ws.Range("G15").Value = 4.1344134611068144
$ws.Range("O15").Value = 3.9902587827274041
$ws.Range("B16").Value = 0.055184677770671972
$ws.Range("C16").Value = 0.055184677770671972
$ws.Range("D16").Value = 2.7335140319094759
$ws.Range("E16").Value = 2.9335140319094761
$ws.Range("F16").Value = 3.5672155534322538
$ws.Range("G16").Value = 3.5672155534322538
$ws.Range("O16").Value = 3.5120308756615821
$ws.Range("B17").Value = 0.059683445724657093
$ws.Range("C17").Value = 0.059683445724657093
$ws.Range("D17").Value = 2.9335140319094761
$ws.Range("E17").Value = 3.1335140319094759
$ws.Range("F17").Value = 3.2670022378070072
$ws.Range("G17").Value = 3.2668728499564779
$ws.Range("O17").Value = 3.2073187920823498
$ws.Range("B18").Value = 0.044684889289771261
$ws.Range("C18").Value = 0.044684889289771261
$ws.Range("D18").Value = 1.433496723860517
$ws.Range("E18").Value = 1.6334967238605169
$ws.Range("F18").Value = 3.3672842405753198
$ws.Range("G18").Value = 3.434199106974714
$ws.Range("O18").Value = 3.322599351285549
$ws.Range("B19").Value = 0.033956332748474462
$ws.Range("C19").Value = 0.033956332748474462
$ws.Range("D19").Value = 0.5334095501305699
$ws.Range("E19").Value = 0.73340955013056985
$ws.Range("F19").Value = 2.0668957364016158
$ws.Range("G19").Value = 2.0668957364016158
$ws.Range("O19").Value = 2.032939403653141
$ws.Range("B20").Value = 0.052501247634525287
$ws.Range("C20").Value = 0.052501247634525287
$ws.Range("D20").Value = 2.0334095501305711
$ws.Range("E20").Value = 2.2334095501305709
$ws.Range("F20").Value = 4.1009008209975351
$ws.Range("G20").Value = 4.1344134611068144
$ws.Range("O20").Value = 4.0483995733630094
$ws.Range("B21").Value = 0.039243527717267607
$ws.Range("C21").Value = 0.039243527717267607
$ws.Range("D21").Value = 0.53349672386051628
$ws.Range("E21").Value = 0.73349672386051623
$ws.Range("F21").Value = 3.5337851608524322
$ws.Range("G21").Value = 3.5672155534322538
$ws.Range("O21").Value = 3.4945416331351651
$ws.Range("B22").Value = 0.054989455760802948
$ws.Range("C22").Value = 0.054989455760802948
$ws.Range("D22").Value = 2.5335140319094762
$ws.Range("E22").Value = 2.7335140319094759
$ws.Range("F22").Value = 3.2668728499564779
$ws.Range("G22").Value = 3.2668728499564779
$ws.Range("O22").Value = 3.211883394195675
$ws.Range("B23").Value = 0.051562136925068308
$ws.Range("C23").Value = 0.051562136925068308
$ws.Range("D23").Value = 1.733496723860517
$ws.Range("E23").Value = 1.933496723860517
$ws.Range("F23").Value = 3.400658243121629
$ws.Range("G23").Value = 3.434199106974714
$ws.Range("O23").Value = 3.3490961061965612
$ws.Range("B24").Value = 0.042478101562498503
$ws.Range("C24").Value = 0.042478101562498503
$ws.Range("D24").Value = 1.1334095501305701
$ws.Range("E24").Value = 1.3334095501305701
$ws.Range("F24").Value = 2.0668327032313871
$ws.Range("G24").Value = 2.0668957364016158
$ws.Range("O24").Value = 2.0243546016688878
$ws.Range("B25").Value = 0.063042567842600969
$ws.Range("C25").Value = 0.063042567842600969
$ws.Range("D25").Value = 3.700345620879903
$ws.Range("E25").Value = 3.9003456208799041
$ws.Range("F25").Value = 4.1344134611068144
$ws.Range("G25").Value = 4.1344134611068144
$ws.Range("O25").Value = 4.0713708932642136
$ws.Range("B26").Value = 0.040282556266250553
$ws.Range("C26").Value = 0.040282556266250553
$ws.Range("D26").Value = 0.83349672386051621
$ws.Range("E26").Value = 1.0334967238605159
$ws.Range("F26").Value = 3.3339197073759972
$ws.Range("G26").Value = 3.434199106974714
$ws.Range("O26").Value = 3.293637151109746
$ws.Range("B27").Value = 0.044508340977005102
$ws.Range("C27").Value = 0.044508340977005102
$ws.Range("D27").Value = 1.4334095501305699
$ws.Range("E27").Value = 1.6334095501305701
$ws.Range("F27").Value = 4.0674217085970774
$ws.Range("G27").Value = 4.1344134611068144
$ws.Range("O27").Value = 4.0229133676200712
$ws.Range("B28").Value = 0.062613184881645273
$ws.Range("C28").Value = 0.062613184881645273
$ws.Range("D28").Value = 2.966914663961731
$ws.Range("E28").Value = 3.1669146639617312
$ws.Range("F28").Value = 3.434199106974714
$ws.Range("G28").Value = 3.434199106974714
$ws.Range("O28").Value = 3.3715859220930691
$ws.Range("B29").Value = 0.040169373417184587
$ws.Range("C29").Value = 0.040169373417184587
$ws.Range("D29").Value = 0.83340955013056983
$ws.Range("E29").Value = 1.03340955013057
$ws.Range("F29").Value = 3.3003751147941678
$ws.Range("G29").Value = 3.434199106974714
$ws.Range("O29").Value = 3.2602057413769829

# These cells no longer have a recorded event in the updated results.
$ws.Range("J3").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("J13").ClearContents()

# Restore the active cell/zoom as left by the editor after the refresh.
$ws.Range("M7").Select()
$excel.ActiveWindow.Zoom = 75
